$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 22; this shifts the existing row 22
# (and everything below it) down by one row, turning the old A22:R134
# block into A23:R135 automatically.
$ws.Rows.Item(22).Insert()

# Populate the newly inserted row 22 with its data. Most fields are the
# same as the row that used to occupy position 22 (now at 23); only the
# date (D), Volumen (J), Precio minimo/maximo/promedio (K/L/M) and
# Precio $/Kg (P) differ.
$ws.Cells.Item(22, 1).Value2 = 5
$ws.Cells.Item(22, 2).Value2 = "Macroferia Regional de Talca"
$ws.Cells.Item(22, 3).Value2 = "Maule"
$ws.Cells.Item(22, 4).Value2 = 45222
$ws.Cells.Item(22, 5).Value2 = 7
$ws.Cells.Item(22, 6).Value2 = 100112026
$ws.Cells.Item(22, 7).Value2 = "Haba"
$ws.Cells.Item(22, 8).Value2 = "Sin especificar"
$ws.Cells.Item(22, 9).Value2 = "Primera"
$ws.Cells.Item(22, 10).Value2 = 500
$ws.Cells.Item(22, 11).Value2 = 10000
$ws.Cells.Item(22, 12).Value2 = 10000
$ws.Cells.Item(22, 13).Value2 = 10000
$ws.Cells.Item(22, 14).Value2 = '$/saco 25 kilos'
$ws.Cells.Item(22, 15).Value2 = "Región del Maule"
$ws.Cells.Item(22, 16).Value2 = 400
$ws.Cells.Item(22, 17).Value2 = 25
$ws.Cells.Item(22, 18).Value2 = "Hortaliza"

# Make sure the date cell keeps the same date number format used by the
# rest of the column (style index 2 in the original workbook).
$ws.Cells.Item(22, 4).NumberFormat = $ws.Cells.Item(23, 4).NumberFormat
